$wb = $excel.ActiveWorkbook

# cs: subjects already correct, just move through it (no cell value changes needed)
$wb.Worksheets.Item("ee").Range("B2").Value = "Linear control"
$wb.Worksheets.Item("ee").Range("B3").Value = "Power systems"
$wb.Worksheets.Item("ee").Range("B4").Select()

$wb.Worksheets.Item("ec").Range("B2").Value = "Communication Systems"
$wb.Worksheets.Item("ec").Range("B3").Value = "Microprocessors and Microcontrollers"
$wb.Worksheets.Item("ec").Range("B4").Select()

$wb.Worksheets.Item("ce").Range("B2").Value = "Geotechnical Engineering"
$wb.Worksheets.Item("ce").Range("B3").Value = "Environmental Engineering"

$wb.Worksheets.Item("me").Range("B2").Value = "Fluid Mechanics"
$wb.Worksheets.Item("me").Range("B3").Value = "Manufacturing Processes"
$wb.Worksheets.Item("me").Range("D3").Select()

$wb.Worksheets.Item("mr").Range("B2").Value = "Mechatronics Design"
$wb.Worksheets.Item("mr").Range("B3").Value = "Industrial Automation"
$wb.Worksheets.Item("mr").Range("D4").Select()

$wb.Worksheets.Item("ad").Range("B2").Value = "Neural Networks and Deep Learning"
$wb.Worksheets.Item("ad").Range("B3").Value = "Natural Language Processing"
$wb.Worksheets.Item("ad").Range("B2:B3").Select()

$wb.Worksheets.Item("rb").Range("B2").Value = "Computer Vision"
$wb.Worksheets.Item("rb").Range("B3").Value = "Robot Perception"
$wb.Worksheets.Item("rb").Range("D5").Select()
